# Server Managment API gegen Brute-Force abgesichert
# Add login_attempts / active columns to the "userx" sheet and a new
# locked-out technical account ("gesperrter_zugang") used to guard the
# server management API against brute-force / dictionary attacks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("userx")

# --- extend header row with the two new columns ---
$ws.Cells.Item(1, 6).Value = "login_attempts"
$ws.Cells.Item(1, 7).Value = "active"

# --- existing rows get default values for the new columns ---
# (values are stored as text, like the rest of this column - column style
# already applies a text number format)
$ws.Cells.Item(2, 6).Value = "0"
$ws.Cells.Item(2, 7).Value = "1"

$ws.Cells.Item(3, 6).Value = "0"
$ws.Cells.Item(3, 7).Value = "1"

# --- new row: a locked technical/guard account ---
$ws.Cells.Item(4, 1).Value = "2"
$ws.Cells.Item(4, 2).Value = "gesperrter_zugang"
$ws.Cells.Item(4, 3).Value = "1"
$ws.Cells.Item(4, 4).Value = "tray_equipment_positionimage"
$ws.Cells.Item(4, 5).Value = "…"
$ws.Cells.Item(4, 6).Value = "5"
$ws.Cells.Item(4, 7).Value = "0"

# --- make userx the active sheet / selection, like in the edited file ---
$ws.Activate()
$ws.Range("C8").Select()
